$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = 'Record'
$ws.Range("B32").Value = 'Balanço Geral'
$ws.Range("C32").Value = 'Limpeza Pública'
$ws.Range("D32").Value = '2025-04-03T13:25'
$ws.Range("E32").Value = 'Negativo'
$ws.Range("F32").Value = 'Imóvel abandonado tomado de mato '

$ws.Range("A33").Value = 'Record'
$ws.Range("B33").Value = 'Balanço Geral'
$ws.Range("C33").Value = 'Limpeza Pública'
$ws.Range("D33").Value = '2025-04-03T13:25'
$ws.Range("E33").Value = 'Negativo'
$ws.Range("F33").Value = 'Imóvel abandonado tomado de mato atrai insetos, ratos e até cobras. Imóvel é na Rua São Jerônimo, no Parque Aurora. Repórter *ao vivo*.  Exibido vídeo de moradora do bairro. Equipe procurou a prefeitura. *sem nota*'

$ws.Range("A34").Value = 'Record'
$ws.Range("B34").Value = 'Balanço Geral'
$ws.Range("C34").Value = 'Esportes'
$ws.Range("D34").Value = '2025-04-03T13:28'
$ws.Range("E34").Value = 'Neutro'
$ws.Range("F34").Value = 'Laila Póvoa assume o Americano com promessa de reestruturar o clube. Sobe som de Laila durante coletiva. *nota coberta*'

